# Update PLC data 2025-10-13 13:50:28
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 162464
$ws.Range("C4").Value = 153467
$ws.Range("C5").Value = 8997
$ws.Range("C8").Value = 64.59
